$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content fixes / additions to the college/department shared-string list
$ws.Range("A4").Value = "CE(Computer Engineering)"
$ws.Range("B5").Value = "MATH(Mathematics)"
$ws.Range("D2").Value = "EARLY(Early Childhood)"
$ws.Range("E3").Value = "ACC(Accounting)"

# A4 ("CIV(Civil Engineering)" -> "CE(Computer Engineering)") also picks up
# a new (black/automatic) font instead of the red "needs review" font used
# by the rest of the data rows.
$ws.Range("A4").Font.Color = 0

# Update the view state: scroll/selection moved from C4 to B2
$excel.Goto($ws.Range("E1"), $true) | Out-Null
$ws.Range("B2").Select() | Out-Null
